# Case_1_29 (380 kV) results: updated res_bus/vm_pu.xlsx values.
# Replaces the flat-start 1.05 pu bus-voltage results (rows 2-25) with the
# converged power-flow results for the 380 kV case.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.055484595288466
$ws.Range("D2").Value = 1.062020755995835
$ws.Range("E2").Value = 1.051825351387883
$ws.Range("F2").Value = 1.070952106439901
$ws.Range("I2").Value = 1.053797396651699
$ws.Range("J2").Value = 1.060491213247869
$ws.Range("K2").Value = 1.064743235087471
$ws.Range("L2").Value = 1.054575780288718
$ws.Range("M2").Value = 1.073650573991002
$ws.Range("N2").Value = 1.061997232410723

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05687599443131
$ws.Range("D3").Value = 1.062876993838295
$ws.Range("E3").Value = 1.053030207113052
$ws.Range("F3").Value = 1.072023716435031
$ws.Range("I3").Value = 1.054237058576753
$ws.Range("J3").Value = 1.061532174587172
$ws.Range("K3").Value = 1.065414103139216
$ws.Range("L3").Value = 1.055592352072423
$ws.Range("M3").Value = 1.074538013597147
$ws.Range("N3").Value = 1.063039672034528

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.057775435110101
$ws.Range("D4").Value = 1.063430483819862
$ws.Range("E4").Value = 1.053809121210794
$ws.Range("F4").Value = 1.072716783528018
$ws.Range("I4").Value = 1.05451987484683
$ws.Range("J4").Value = 1.06220439953324
$ws.Range("K4").Value = 1.065847003081147
$ws.Range("L4").Value = 1.056248867329686
$ws.Range("M4").Value = 1.075111308902807
$ws.Range("N4").Value = 1.063712851617125

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.058153352245308
$ws.Range("D5").Value = 1.063663039783143
$ws.Range("E5").Value = 1.054136410567349
$ws.Range("F5").Value = 1.073008070321341
$ws.Range("I5").Value = 1.05463837126702
$ws.Range("J5").Value = 1.062486683786177
$ws.Range("K5").Value = 1.066028709117001
$ws.Range("L5").Value = 1.056524564304426
$ws.Range("M5").Value = 1.0753520996493
$ws.Range("N5").Value = 1.063995536746079

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.058216794152391
$ws.Range("D6").Value = 1.063702079292328
$ws.Range("E6").Value = 1.054191354246988
$ws.Range("F6").Value = 1.073056974100274
$ws.Range("I6").Value = 1.054658243946246
$ws.Range("J6").Value = 1.062534061936313
$ws.Range("K6").Value = 1.066059201680843
$ws.Range("L6").Value = 1.056570837418278
$ws.Range("M6").Value = 1.075392516474797
$ws.Range("N6").Value = 1.064042982178623

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.057780485667654
$ws.Range("D7").Value = 1.063433591758169
$ws.Range("E7").Value = 1.053813495115757
$ws.Range("F7").Value = 1.072720676022971
$ws.Range("I7").Value = 1.05452145977009
$ws.Range("J7").Value = 1.062208172679085
$ws.Range("K7").Value = 1.065849432164664
$ws.Range("L7").Value = 1.056252552389099
$ws.Range("M7").Value = 1.075114527232949
$ws.Range("N7").Value = 1.06371663012127

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.055955011802614
$ws.Range("D8").Value = 1.062310240845352
$ws.Range("E8").Value = 1.05223268669033
$ws.Range("F8").Value = 1.071314332809851
$ws.Range("I8").Value = 1.0539463299449
$ws.Range("J8").Value = 1.060843291813782
$ws.Range("K8").Value = 1.064970206658431
$ws.Range("L8").Value = 1.054919601289604
$ws.Range("M8").Value = 1.073950683573067
$ws.Range("N8").Value = 1.062349810968584

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.052731281675148
$ws.Range("D9").Value = 1.060326461861782
$ws.Range("E9").Value = 1.049441524334526
$ws.Range("F9").Value = 1.068833516235759
$ws.Range("I9").Value = 1.052919994641413
$ws.Range("J9").Value = 1.05842773644833
$ws.Range("K9").Value = 1.063411673547697
$ws.Range("L9").Value = 1.052560868701486
$ws.Range("M9").Value = 1.071892579896774
$ws.Range("N9").Value = 1.059930825237445

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.050577112957883
$ws.Range("D10").Value = 1.059000978296841
$ws.Range("E10").Value = 1.04757678510104
$ws.Range("F10").Value = 1.067177728297975
$ws.Range("I10").Value = 1.052227025063215
$ws.Range("J10").Value = 1.056810125164229
$ws.Range("K10").Value = 1.06236636098221
$ws.Range("L10").Value = 1.050981524784366
$ws.Range("M10").Value = 1.070515509978622
$ws.Range("N10").Value = 1.058310916759778

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.049643077074695
$ws.Range("D11").Value = 1.058426307747034
$ws.Range("E11").Value = 1.046768343570363
$ws.Range("F11").Value = 1.066460272818242
$ws.Range("I11").Value = 1.051924867954077
$ws.Range("J11").Value = 1.056107917315635
$ws.Range("K11").Value = 1.061912216876012
$ws.Range("L11").Value = 1.050295980822018
$ws.Range("M11").Value = 1.069918010379036
$ws.Range("N11").Value = 1.057607711695494

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.049295938368985
$ws.Range("D12").Value = 1.058212738542805
$ws.Range("E12").Value = 1.046467898309253
$ws.Range("F12").Value = 1.066193702026845
$ws.Range("I12").Value = 1.051812316619912
$ws.Range("J12").Value = 1.055846815785593
$ws.Range("K12").Value = 1.06174329781446
$ws.Range("L12").Value = 1.050041083299969
$ws.Range("M12").Value = 1.069695886825771
$ws.Range("N12").Value = 1.057346239371331

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.049370409772786
$ws.Range("D13").Value = 1.058258554911214
$ws.Range("E13").Value = 1.046532351878967
$ws.Range("F13").Value = 1.066250885840775
$ws.Range("I13").Value = 1.051836473623649
$ws.Range("J13").Value = 1.055902835246813
$ws.Range("K13").Value = 1.061779541951362
$ws.Range("L13").Value = 1.050095771330192
$ws.Range("M13").Value = 1.069743541519742
$ws.Range("N13").Value = 1.057402338386612

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.049614386517433
$ws.Range("D14").Value = 1.058408656323665
$ws.Range("E14").Value = 1.046743511839348
$ws.Range("F14").Value = 1.066438239556005
$ws.Range("I14").Value = 1.051915570896891
$ws.Range("J14").Value = 1.056086340119051
$ws.Range("K14").Value = 1.061898258670521
$ws.Range("L14").Value = 1.050274916157837
$ws.Range("M14").Value = 1.069899653377048
$ws.Range("N14").Value = 1.057586103856816

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.049764682467993
$ws.Range("D15").Value = 1.058501124001831
$ws.Range("E15").Value = 1.046873593914714
$ws.Range("F15").Value = 1.066553664233015
$ws.Range("I15").Value = 1.051964263315911
$ws.Range("J15").Value = 1.056199367590237
$ws.Range("K15").Value = 1.061971373464349
$ws.Range("L15").Value = 1.05038525916228
$ws.Range("M15").Value = 1.069995814378369
$ws.Range("N15").Value = 1.057699291839976

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.05063907465457
$ws.Range("D16").Value = 1.059039101818969
$ws.Range("E16").Value = 1.047630417330226
$ws.Range("F16").Value = 1.06722533295207
$ws.Range("I16").Value = 1.052247033913294
$ws.Range("J16").Value = 1.056856690711827
$ws.Range("K16").Value = 1.06239646891651
$ws.Range("L16").Value = 1.051026986399066
$ws.Range("M16").Value = 1.07055513816498
$ws.Range("N16").Value = 1.058357548435794

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.051187214509654
$ws.Range("D17").Value = 1.059376365152306
$ws.Range("E17").Value = 1.048104882692404
$ws.Range("F17").Value = 1.067646520476502
$ws.Range("I17").Value = 1.052423845875983
$ws.Range("J17").Value = 1.057268534957721
$ws.Range("K17").Value = 1.062662712667902
$ws.Range("L17").Value = 1.051429073026769
$ws.Range("M17").Value = 1.070905659069875
$ws.Range("N17").Value = 1.058769977547757

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.051506813404152
$ws.Range("D18").Value = 1.059573015079052
$ws.Range("E18").Value = 1.048381534392224
$ws.Range("F18").Value = 1.067892145083862
$ws.Range("I18").Value = 1.052526775117859
$ws.Range("J18").Value = 1.057508586096068
$ws.Range("K18").Value = 1.062817861912821
$ws.Range("L18").Value = 1.051663441754618
$ws.Range("M18").Value = 1.071109994509973
$ws.Range("N18").Value = 1.059010369586255

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.051615767884617
$ws.Range("D19").Value = 1.059640055787617
$ws.Range("E19").Value = 1.048475849331597
$ws.Range("F19").Value = 1.067975888839153
$ws.Range("I19").Value = 1.052561837087176
$ws.Range("J19").Value = 1.057590408510908
$ws.Range("K19").Value = 1.062870739047422
$ws.Range("L19").Value = 1.051743328142273
$ws.Range("M19").Value = 1.071179647761656
$ws.Range("N19").Value = 1.059092308198309

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.051128416934324
$ws.Range("D20").Value = 1.059340187249156
$ws.Range("E20").Value = 1.048053987002899
$ws.Range("F20").Value = 1.067601335916653
$ws.Range("I20").Value = 1.05240489655171
$ws.Range("J20").Value = 1.057224365649229
$ws.Range("K20").Value = 1.062634162372416
$ws.Range("L20").Value = 1.051385949681531
$ws.Range("M20").Value = 1.070868063663388
$ws.Range("N20").Value = 1.05872574551378

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.049542546910422
$ws.Range("D21").Value = 1.058364458275963
$ws.Range("E21").Value = 1.046681334781847
$ws.Range("F21").Value = 1.066383070671478
$ws.Range("I21").Value = 1.051892287478566
$ws.Range("J21").Value = 1.056032310001895
$ws.Range("K21").Value = 1.061863305911798
$ws.Range("L21").Value = 1.050222169568841
$ws.Range("M21").Value = 1.069853687469011
$ws.Range("N21").Value = 1.057531997010697

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.048544309443231
$ws.Range("D22").Value = 1.057750335830301
$ws.Range("E22").Value = 1.045817400254346
$ws.Range("F22").Value = 1.065616658676261
$ws.Range("I22").Value = 1.051568156444444
$ws.Range("J22").Value = 1.055281251551425
$ws.Range("K22").Value = 1.061377307857607
$ws.Range("L22").Value = 1.049488971895735
$ws.Range("M22").Value = 1.069214833869165
$ws.Range("N22").Value = 1.056779871971079

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04907360386037
$ws.Range("D23").Value = 1.05807595526637
$ws.Range("E23").Value = 1.046275474427649
$ws.Range("F23").Value = 1.066022990751364
$ws.Range("I23").Value = 1.05174015879783
$ws.Range("J23").Value = 1.055679551643482
$ws.Range("K23").Value = 1.061635071303521
$ws.Range("L23").Value = 1.049877795734648
$ws.Range("M23").Value = 1.069553604962324
$ws.Range("N23").Value = 1.057178737694953

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.05115498540371
$ws.Range("D24").Value = 1.059356534702852
$ws.Range("E24").Value = 1.048076984870367
$ws.Range("F24").Value = 1.067621753019338
$ws.Range("I24").Value = 1.052413459560331
$ws.Range("J24").Value = 1.05724432438552
$ws.Range("K24").Value = 1.062647063473809
$ws.Range("L24").Value = 1.051405435764128
$ws.Range("M24").Value = 1.07088505177222
$ws.Range("N24").Value = 1.058745732593767

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.053565555029832
$ws.Range("D25").Value = 1.060839832896591
$ws.Range("E25").Value = 1.050163789593386
$ws.Range("F25").Value = 1.069475194496039
$ws.Range("I25").Value = 1.053186862265363
$ws.Range("J25").Value = 1.059053476401751
$ws.Range("K25").Value = 1.063815693976133
$ws.Range("L25").Value = 1.053171851397408
$ws.Range("M25").Value = 1.072425521756394
$ws.Range("N25").Value = 1.060557453813373
